# Generate Report for Handback
# Update the timestamp values on the Overview / zh-cn / de-de sheets to
# reflect a newer handback report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for 688fd49d-...md (row 3, column G)
$wsOverview.Range("G3").Value = "2016-08-19 14:51:58"

# de-de: "Correspond Handoff Datetime" for 688fd49d-...md (row 3, column H)
$wsDeDe.Range("H3").Value = "2016-08-19 14:51:58"

# zh-cn: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for 688fd49d-...md (row 3)
$wsZhCn.Range("H3").Value = "2016-08-19 14:51:54"
$wsZhCn.Range("K3").Value = "2016-08-19 14:52:29"

# de-de: "Correspond Handback DateTime" for 688fd49d-...md (row 3, column K)
$wsDeDe.Range("K3").Value = "2016-08-19 14:52:36"
